# chore: update Sheets via scheduled runner
# Refreshes cached market-price / leve-profit figures (columns H-N) on a
# handful of rows across the ALC/BSM/CRP/CUL/GSM/LTW/WVR leve tables.
# Some rows collapse from 7 data columns (H..N) down to 6 because the
# HQ-profit figure (col N, or col M on rows that already lacked an N
# value) no longer applies - those trailing cells are cleared outright
# rather than zeroed, matching upstream's column removal.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 43
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()

# Row 58
$ws.Range("H58").Value = 3584.875
$ws.Range("I58").Value = 286.25
$ws.Range("J58").Value = 6883.5
$ws.Range("K58").Value = 858.75
$ws.Range("L58").Value = 20650.5
$ws.Range("M58").Value = -708.75
$ws.Range("N58").Value = -20950.5

# Row 76
$ws.Range("H76").Value = 3192.4666
$ws.Range("I76").Value = 3129
$ws.Range("J76").Value = 3208.3333
$ws.Range("K76").Value = 3129
$ws.Range("L76").Value = 3208.3333
$ws.Range("M76").Value = -2814
$ws.Range("N76").Value = -3838.3333

# Row 79
$ws.Range("H79").Value = 3192.4666
$ws.Range("I79").Value = 3129
$ws.Range("J79").Value = 3208.3333
$ws.Range("K79").Value = 3129
$ws.Range("L79").Value = 3208.3333
$ws.Range("M79").Value = -2037
$ws.Range("N79").Value = -5392.3333

# Row 132
$ws.Range("H132").Value = 4000.9092
$ws.Range("I132").Value = 4223
$ws.Range("K132").Value = 12669
$ws.Range("M132").Value = -10139

# Row 137
$ws.Range("H137").Value = 115122.336
$ws.Range("I137").Value = 8000
$ws.Range("J137").Value = 168683.5
$ws.Range("K137").Value = 24000
$ws.Range("L137").Value = 506050.5
$ws.Range("M137").Value = -21450
$ws.Range("N137").Value = -511150.5


# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 20343.893
$ws.Range("I134").Value = 27323.365
$ws.Range("K134").Value = 81970.095
$ws.Range("M134").Value = -79435.095


# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 80
$ws.Range("H80").Value = 20509.334
$ws.Range("J80").Value = 20509.334
$ws.Range("L80").Value = 20509.334
$ws.Range("N80").Value = -22755.334

# Row 83
$ws.Range("H83").Value = 20509.334
$ws.Range("J83").Value = 20509.334
$ws.Range("L83").Value = 61528.00199999999
$ws.Range("N83").Value = -72760.00199999999

# Row 99
$ws.Range("H99").Value = 13892676
$ws.Range("I99").Value = 3073.16
$ws.Range("J99").Value = 45459950
$ws.Range("K99").Value = 3073.16
$ws.Range("L99").Value = 45459950
$ws.Range("M99").Value = -1575.16
$ws.Range("N99").Value = -45462946

# Row 126
$ws.Range("H126").Value = 13892676
$ws.Range("I126").Value = 3073.16
$ws.Range("J126").Value = 45459950
$ws.Range("K126").Value = 9219.48
$ws.Range("L126").Value = 136379850
$ws.Range("M126").Value = -6749.48
$ws.Range("N126").Value = -136384790

# Row 132
$ws.Range("H132").Value = 19095.867
$ws.Range("I132").Value = 27723.21
$ws.Range("K132").Value = 83169.63
$ws.Range("M132").Value = -80639.63


# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 655.7273
$ws.Range("I5").Value = 655.7273
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 1967.1819
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -1855.1819
$ws.Range("N5").ClearContents()

# Row 131
$ws.Range("H131").Value = 737.77
$ws.Range("J131").Value = 737.77
$ws.Range("L131").Value = 2213.31
$ws.Range("N131").Value = -12293.31

# Row 135
$ws.Range("H135").Value = 655.7273
$ws.Range("I135").Value = 655.7273
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 5901.545700000001
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -3366.545700000001
$ws.Range("N135").ClearContents()


# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 15
$ws.Range("H15").Value = 16000
$ws.Range("J15").Value = 16000
$ws.Range("L15").Value = 16000
$ws.Range("N15").Value = -16576

# Row 70
$ws.Range("H70").Value = 11765.154
$ws.Range("I70").Value = 19899.5
$ws.Range("J70").Value = 4792.857
$ws.Range("K70").Value = 19899.5
$ws.Range("L70").Value = 4792.857
$ws.Range("M70").Value = -19629.5
$ws.Range("N70").Value = -5332.857

# Row 73
$ws.Range("H73").Value = 11765.154
$ws.Range("I73").Value = 19899.5
$ws.Range("J73").Value = 4792.857
$ws.Range("K73").Value = 19899.5
$ws.Range("L73").Value = 4792.857
$ws.Range("M73").Value = -18963.5
$ws.Range("N73").Value = -6664.857

# Row 80
$ws.Range("H80").Value = 3999.5881
$ws.Range("I80").Value = 3825
$ws.Range("J80").Value = 4053.3076
$ws.Range("K80").Value = 3825
$ws.Range("L80").Value = 4053.3076
$ws.Range("M80").Value = -2827
$ws.Range("N80").Value = -6049.3076

# Row 81
$ws.Range("H81").Value = 16000
$ws.Range("J81").Value = 16000
$ws.Range("L81").Value = 16000
$ws.Range("N81").Value = -17996

# Row 83
$ws.Range("H83").Value = 3999.5881
$ws.Range("I83").Value = 3825
$ws.Range("J83").Value = 4053.3076
$ws.Range("K83").Value = 19125
$ws.Range("L83").Value = 20266.538
$ws.Range("M83").Value = -14133
$ws.Range("N83").Value = -30250.538

# Row 84
$ws.Range("H84").Value = 16000
$ws.Range("J84").Value = 16000
$ws.Range("L84").Value = 48000
$ws.Range("N84").Value = -57984

# Row 132
$ws.Range("H132").Value = 94328.47
$ws.Range("I132").Value = 90194.664
$ws.Range("K132").Value = 270583.992
$ws.Range("M132").Value = -268053.992


# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 4586.4443
$ws.Range("J40").Value = 4828.7144
$ws.Range("L40").Value = 4828.7144
$ws.Range("N40").Value = -5100.7144

# Row 46
$ws.Range("H46").Value = 1122.4706
$ws.Range("I46").Value = 914
$ws.Range("J46").Value = 1800
$ws.Range("K46").Value = 914
$ws.Range("L46").Value = 1800
$ws.Range("M46").Value = -726
$ws.Range("N46").Value = -2176

# Row 61
$ws.Range("H61").Value = 4203.143
$ws.Range("I61").Value = 1947
$ws.Range("J61").Value = 7211.3335
$ws.Range("K61").Value = 1947
$ws.Range("L61").Value = 7211.3335
$ws.Range("M61").Value = -1745
$ws.Range("N61").Value = -7615.3335

# Row 113
$ws.Range("H113").Value = 4203.143
$ws.Range("I113").Value = 1947
$ws.Range("J113").Value = 7211.3335
$ws.Range("K113").Value = 1947
$ws.Range("L113").Value = 7211.3335
$ws.Range("M113").Value = 223
$ws.Range("N113").Value = -11551.3335

# Row 132
$ws.Range("H132").Value = 1938.6666
$ws.Range("I132").Value = 1206.7142
$ws.Range("K132").Value = 3620.1426
$ws.Range("M132").Value = -1090.1426


# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 75
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()

# Row 78
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

# Row 122
$ws.Range("H122").Value = 1555.8572
$ws.Range("I122").Value = 1662.8695
$ws.Range("K122").Value = 4988.6085
$ws.Range("M122").Value = -2538.6085

